$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the edited cells keep their text/string representation
$ws.Range("C2:F5").NumberFormat = "@"

# Row 2
$ws.Range("C2").Value = "30"
$ws.Range("D2").Value = "31"
$ws.Range("E2").Value = "1"
$ws.Range("F2").Value = "1"

# Row 3
$ws.Range("C3").Value = "58"
$ws.Range("D3").Value = "45"
$ws.Range("E3").Value = "7"
$ws.Range("F3").Value = "1"

# Row 4
$ws.Range("C4").Value = "39"
$ws.Range("D4").Value = "32"
$ws.Range("E4").Value = "4"
$ws.Range("F4").Value = "1"

# Row 5
$ws.Range("C5").Value = "87"
$ws.Range("D5").Value = "45"
$ws.Range("E5").Value = "12"
$ws.Range("F5").Value = "2"
